# Updated cryptos list on Fri Apr 21 19:22:46 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns with the latest scraped
# values. A couple of coins swapped rank positions (B/C columns updated
# accordingly). Numeric-looking price strings are prefixed with a leading
# apostrophe so Excel stores them as text (matching the original inline
# string cells) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.572.48'
$ws.Range('E2').Value = '  -3.46%  '
$ws.Range('D3').Value = '1.861.01'
$ws.Range('E3').Value = '  -4.73%  '
$ws.Range('D4').Value = '''1.008'
$ws.Range('E4').Value = '  -0.70%  '
$ws.Range('D5').Value = '''323.32'
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('D6').Value = '''1.007'
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('D7').Value = '''0.4494'
$ws.Range('E7').Value = '  -5.61%  '
$ws.Range('D8').Value = '''0.3856'
$ws.Range('E8').Value = '  -4.40%  '
$ws.Range('D9').Value = '''48.53'
$ws.Range('E9').Value = '  -9.69%  '
$ws.Range('D10').Value = '''0.08028'
$ws.Range('E10').Value = '  -5.42%  '
$ws.Range('D11').Value = '''1.018'
$ws.Range('E11').Value = '  -3.88%  '
$ws.Range('E12').Value = '  -2.85%  '
$ws.Range('D13').Value = '1.879.40'
$ws.Range('E13').Value = '  -4.11%  '
$ws.Range('D14').Value = '''7.181'
$ws.Range('E14').Value = '  -5.53%  '
$ws.Range('D15').Value = '''5.873'
$ws.Range('E15').Value = '  -5.36%  '
$ws.Range('D16').Value = '''1.010'
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '''0.00001037'
$ws.Range('E17').Value = '  -3.72%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').Value = '''86.12'
$ws.Range('E18').Value = '  -3.16%  '
$ws.Range('D19').Value = '''0.06546'
$ws.Range('E19').Value = '  -1.06%  '
$ws.Range('D20').Value = '''17.13'
$ws.Range('E20').Value = '  -8.22%  '
$ws.Range('E21').Value = '  -0.53%  '
$ws.Range('D22').Value = '''5.518'
$ws.Range('E22').Value = '  -5.09%  '
$ws.Range('D23').Value = '27.606.19'
$ws.Range('E23').Value = '  -3.44%  '
$ws.Range('D24').Value = '''10.86'
$ws.Range('E24').Value = '  -5.61%  '
$ws.Range('D25').Value = '''2.306'
$ws.Range('E25').Value = '  +0.30%  '
$ws.Range('D26').Value = '2.120.91'
$ws.Range('E26').Value = '  -3.23%  '
$ws.Range('D27').Value = '''151.27'
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('D28').Value = '''19.44'
$ws.Range('E28').Value = '  -3.32%  '
$ws.Range('D29').Value = '''5.531'
$ws.Range('E29').Value = '  -6.57%  '
$ws.Range('D30').Value = '''2.031'
$ws.Range('E30').Value = '  -5.76%  '
$ws.Range('D31').Value = '''120.06'
$ws.Range('E31').Value = '  -2.83%  '
$ws.Range('D32').Value = '''0.09390'
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('E33').Value = '  +1.46%  '
$ws.Range('D34').Value = '''0.9270'
$ws.Range('E34').Value = '  -7.04%  '
$ws.Range('D35').Value = '''3.630'
$ws.Range('E35').Value = '  -1.14%  '
$ws.Range('D36').Value = '''5.274'
$ws.Range('E36').Value = '  -5.80%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.02232'
$ws.Range('E37').Value = '  -4.20%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '''1.222'
$ws.Range('E38').Value = '  -2.83%  '
$ws.Range('D39').Value = '''0.05977'
$ws.Range('E39').Value = '  -3.97%  '
$ws.Range('D40').Value = '''8.397'
$ws.Range('E40').Value = '  -3.71%  '
$ws.Range('D41').Value = '''1.008'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('D42').Value = '''0.5943'
$ws.Range('E42').Value = '  -4.58%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '''0.1854'
$ws.Range('E43').Value = '  -3.59%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = '''10.30'
$ws.Range('E44').Value = '  -7.09%  '
$ws.Range('D45').Value = '''1.278'
$ws.Range('E45').Value = '  -3.89%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.5652'
$ws.Range('E46').Value = '  -4.91%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''12.27'
$ws.Range('E47').Value = '  -5.10%  '
$ws.Range('D48').Value = '''3.408'
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('D49').Value = '''1.938'
$ws.Range('E49').Value = '  -6.11%  '
$ws.Range('D50').Value = '''0.06865'
$ws.Range('E50').Value = '  +0.70%  '
$ws.Range('E51').Value = '  -0.74%  '
